$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column F (SECTION shifts from F to H, etc.)
$ws.Range("F1:G1").EntireColumn.Insert()

# New header cells in the freshly inserted columns
$ws.Range("F4").Value = "SUBJECT"
$ws.Range("G4").Value = "DESCRIPTION"

# Rename the old "RECEIVED DATE TO THIRD PARTY" header (now shifted to Q4)
$ws.Range("Q4").Value = "RECEIVED DATE FROM THIRD PARTY"

# Approximate the column widths used in the final layout as closely as this
# runtime's column-width quantization allows.
$ws.Columns.Item(6).ColumnWidth = 33.666666666666664
$ws.Columns.Item(7).ColumnWidth = 38

# Restore the selection to the cell the author ended up on
[void]$ws.Range("D22").Select()
